# Auto-generated Excel COM-interop edit script
# Applies the data refresh described by the diff to sheet '展览' (index 1) and '全部类型' (index 4)

function Set-TextCell {
    param($ws, [string]$addr, [string]$val)
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

function Set-NumCell {
    param($ws, [string]$addr, $val)
    $ws.Range($addr).Value = $val
}

$wb = $excel.ActiveWorkbook
$wsExhibit = $wb.Worksheets.Item(1)   # 展览
$wsAll     = $wb.Worksheets.Item(4)   # 全部类型

# ---- sheet1 (展览) rows 27-31 full content refresh ----
# Row 27: 杭州·2024ESCC游戏电竞博览会暨新次元微光青春动漫交流会
Set-TextCell $wsExhibit "B27" '2024-04-04'
Set-TextCell $wsExhibit "C27" '杭州·2024ESCC游戏电竞博览会暨新次元微光青春动漫交流会'
Set-TextCell $wsExhibit "D27" '钱江世纪城奔竞大道353号 杭州国际博览中心'
Set-TextCell $wsExhibit "E27" '2024.04.04 09:30-04.05 16:30'
Set-NumCell $wsExhibit "F27" 1466
Set-NumCell $wsExhibit "G27" 75
Set-TextCell $wsExhibit "H27" 'https://show.bilibili.com/platform/detail.html?id=81450'
Set-TextCell $wsExhibit "I27" '//i1.hdslb.com/bfs/openplatform/202403/OfpkJ50P1709548942017.png'

# Row 28: 杭州·ELECTRIC COMIC动漫游戏展
Set-TextCell $wsExhibit "B28" '2024-04-04'
Set-TextCell $wsExhibit "C28" '杭州·ELECTRIC COMIC动漫游戏展'
Set-TextCell $wsExhibit "D28" '沈半路171号 T-Car杭州汽车文化主题公园'
Set-TextCell $wsExhibit "E28" '2024.04.04 10:00-04.05 17:00'
Set-NumCell $wsExhibit "F28" 1
Set-NumCell $wsExhibit "G28" 70
Set-TextCell $wsExhibit "H28" 'https://show.bilibili.com/platform/detail.html?id=82270'
Set-TextCell $wsExhibit "I28" '//i1.hdslb.com/bfs/openplatform/202403/Be5zFgv11709543746638.jpeg'

# Row 29: 杭州·创造力动漫游戏嘉年华1.0
Set-TextCell $wsExhibit "B29" '2024-04-04'
Set-TextCell $wsExhibit "C29" '杭州·创造力动漫游戏嘉年华1.0'
Set-TextCell $wsExhibit "D29" '沈半路171号 T-Car杭州汽车文化主题公园'
Set-TextCell $wsExhibit "E29" '2024.04.04 10:00-04.05 17:00'
Set-NumCell $wsExhibit "F29" 63
Set-TextCell $wsExhibit "G29" '不可售'
Set-TextCell $wsExhibit "H29" 'https://show.bilibili.com/platform/detail.html?id=81078'
Set-TextCell $wsExhibit "I29" '//i0.hdslb.com/bfs/openplatform/202401/o4cl1vwE1705635692432.jpeg'

# Row 30: 杭州·梦漫星河动漫展
Set-TextCell $wsExhibit "B30" '2024-04-04'
Set-TextCell $wsExhibit "C30" '杭州·梦漫星河动漫展'
Set-TextCell $wsExhibit "D30" '德胜东路2539号 梦马汽车小镇'
Set-TextCell $wsExhibit "E30" '2024.04.04 10:00-04.05 17:00'
Set-NumCell $wsExhibit "F30" 1657
Set-NumCell $wsExhibit "G30" 58.5
Set-TextCell $wsExhibit "H30" 'https://show.bilibili.com/platform/detail.html?id=81699'
Set-TextCell $wsExhibit "I30" '//i0.hdslb.com/bfs/openplatform/202402/sZfZd47Y1706868453434.jpeg'

# Row 31: 杭州·第九届萌次元动漫嘉年华
Set-TextCell $wsExhibit "B31" '2024-04-04'
Set-TextCell $wsExhibit "C31" '杭州·第九届萌次元动漫嘉年华'
Set-TextCell $wsExhibit "D31" '长乐路29号五组2幢 杭州运河文化发布中心'
Set-TextCell $wsExhibit "E31" '2024.04.04 10:00-04.05 17:00'
Set-NumCell $wsExhibit "F31" 240
Set-TextCell $wsExhibit "G31" '不可售'
Set-TextCell $wsExhibit "H31" 'https://show.bilibili.com/platform/detail.html?id=78866'
Set-TextCell $wsExhibit "I31" '//i1.hdslb.com/bfs/openplatform/202311/8jSeAOZH1700636327971.jpeg'

# ---- sheet4 (全部类型) rows 31-35 full content refresh ----
# Row 31: 杭州·2024ESCC游戏电竞博览会暨新次元微光青春动漫交流会
Set-TextCell $wsAll "B31" '2024-04-04'
Set-TextCell $wsAll "C31" '杭州·2024ESCC游戏电竞博览会暨新次元微光青春动漫交流会'
Set-TextCell $wsAll "D31" '钱江世纪城奔竞大道353号 杭州国际博览中心'
Set-TextCell $wsAll "E31" '2024.04.04 09:30-04.05 16:30'
Set-NumCell $wsAll "F31" 1466
Set-NumCell $wsAll "G31" 75
Set-TextCell $wsAll "H31" 'https://show.bilibili.com/platform/detail.html?id=81450'
Set-TextCell $wsAll "I31" '//i1.hdslb.com/bfs/openplatform/202403/OfpkJ50P1709548942017.png'

# Row 32: 杭州·ELECTRIC COMIC动漫游戏展
Set-TextCell $wsAll "B32" '2024-04-04'
Set-TextCell $wsAll "C32" '杭州·ELECTRIC COMIC动漫游戏展'
Set-TextCell $wsAll "D32" '沈半路171号 T-Car杭州汽车文化主题公园'
Set-TextCell $wsAll "E32" '2024.04.04 10:00-04.05 17:00'
Set-NumCell $wsAll "F32" 1
Set-NumCell $wsAll "G32" 70
Set-TextCell $wsAll "H32" 'https://show.bilibili.com/platform/detail.html?id=82270'
Set-TextCell $wsAll "I32" '//i1.hdslb.com/bfs/openplatform/202403/Be5zFgv11709543746638.jpeg'

# Row 33: 杭州·创造力动漫游戏嘉年华1.0
Set-TextCell $wsAll "B33" '2024-04-04'
Set-TextCell $wsAll "C33" '杭州·创造力动漫游戏嘉年华1.0'
Set-TextCell $wsAll "D33" '沈半路171号 T-Car杭州汽车文化主题公园'
Set-TextCell $wsAll "E33" '2024.04.04 10:00-04.05 17:00'
Set-NumCell $wsAll "F33" 63
Set-TextCell $wsAll "G33" '不可售'
Set-TextCell $wsAll "H33" 'https://show.bilibili.com/platform/detail.html?id=81078'
Set-TextCell $wsAll "I33" '//i0.hdslb.com/bfs/openplatform/202401/o4cl1vwE1705635692432.jpeg'

# Row 34: 杭州·梦漫星河动漫展
Set-TextCell $wsAll "B34" '2024-04-04'
Set-TextCell $wsAll "C34" '杭州·梦漫星河动漫展'
Set-TextCell $wsAll "D34" '德胜东路2539号 梦马汽车小镇'
Set-TextCell $wsAll "E34" '2024.04.04 10:00-04.05 17:00'
Set-NumCell $wsAll "F34" 1657
Set-NumCell $wsAll "G34" 58.5
Set-TextCell $wsAll "H34" 'https://show.bilibili.com/platform/detail.html?id=81699'
Set-TextCell $wsAll "I34" '//i0.hdslb.com/bfs/openplatform/202402/sZfZd47Y1706868453434.jpeg'

# Row 35: 杭州·第九届萌次元动漫嘉年华
Set-TextCell $wsAll "B35" '2024-04-04'
Set-TextCell $wsAll "C35" '杭州·第九届萌次元动漫嘉年华'
Set-TextCell $wsAll "D35" '长乐路29号五组2幢 杭州运河文化发布中心'
Set-TextCell $wsAll "E35" '2024.04.04 10:00-04.05 17:00'
Set-NumCell $wsAll "F35" 240
Set-TextCell $wsAll "G35" '不可售'
Set-TextCell $wsAll "H35" 'https://show.bilibili.com/platform/detail.html?id=78866'
Set-TextCell $wsAll "I35" '//i1.hdslb.com/bfs/openplatform/202311/8jSeAOZH1700636327971.jpeg'

# ---- simple '想去人数' (F column) increments: sheet1 (展览) ----
Set-NumCell $wsExhibit "F3" 3167   # was 3162
Set-NumCell $wsExhibit "F5" 2201   # was 2193
Set-NumCell $wsExhibit "F6" 329   # was 328
Set-NumCell $wsExhibit "F8" 1049   # was 1044
Set-NumCell $wsExhibit "F9" 1016   # was 1013
Set-NumCell $wsExhibit "F11" 463   # was 460
Set-NumCell $wsExhibit "F12" 1158   # was 1156
Set-NumCell $wsExhibit "F16" 7834   # was 7811
Set-NumCell $wsExhibit "F17" 342   # was 341
Set-NumCell $wsExhibit "F18" 2466   # was 2463
Set-NumCell $wsExhibit "F19" 214   # was 213
Set-NumCell $wsExhibit "F25" 1133   # was 1131
Set-NumCell $wsExhibit "F33" 482   # was 480
Set-NumCell $wsExhibit "F36" 272   # was 271
Set-NumCell $wsExhibit "F37" 42   # was 41
Set-NumCell $wsExhibit "F38" 172   # was 169
Set-NumCell $wsExhibit "F39" 343   # was 342
Set-NumCell $wsExhibit "F41" 218   # was 217

# ---- simple '想去人数' (F column) increments: sheet4 (全部类型) ----
Set-NumCell $wsAll "F5" 3167   # was 3162
Set-NumCell $wsAll "F7" 2201   # was 2193
Set-NumCell $wsAll "F8" 329   # was 328
Set-NumCell $wsAll "F10" 1049   # was 1044
Set-NumCell $wsAll "F12" 1016   # was 1013
Set-NumCell $wsAll "F14" 463   # was 460
Set-NumCell $wsAll "F15" 1158   # was 1156
Set-NumCell $wsAll "F19" 7834   # was 7812
Set-NumCell $wsAll "F20" 342   # was 341
Set-NumCell $wsAll "F21" 2466   # was 2463
Set-NumCell $wsAll "F23" 214   # was 213
Set-NumCell $wsAll "F29" 1133   # was 1131
Set-NumCell $wsAll "F37" 482   # was 480
Set-NumCell $wsAll "F40" 272   # was 271
Set-NumCell $wsAll "F41" 42   # was 41
Set-NumCell $wsAll "F42" 172   # was 169
Set-NumCell $wsAll "F43" 343   # was 342
Set-NumCell $wsAll "F48" 218   # was 217

Write-Host "Edit applied successfully"